# Se añade el metodo /get-areas
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Dates in column A are stored as plain text (not date serials) in the source
# data, so mark the column as Text before writing date-like strings into it.
$ws.Range("A2:A7").NumberFormat = "@"

# Row 2
$ws.Cells.Item(2, 1).Value = '2024-08-28'
$ws.Cells.Item(2, 2).Value = 'II'
$ws.Cells.Item(2, 3).Value = 'MOLINO'
$ws.Cells.Item(2, 4).Value = 'BARRE SARANGO JONATHAN VINICIO'
$ws.Cells.Item(2, 5).Value = 'Cumple'
$ws.Cells.Item(2, 6).Value = 'Cumple'
$ws.Cells.Item(2, 7).Value = 'Cumple'
$ws.Cells.Item(2, 8).Value = 'No cumple'
$ws.Cells.Item(2, 9).Value = 'Cumple'
$ws.Cells.Item(2, 10).Value = 'Cumple'
$ws.Cells.Item(2, 11).Value = 'Cumple'
$ws.Cells.Item(2, 12).Value = 'Cumple'
$ws.Cells.Item(2, 13).Value = 'No cumple'
$ws.Cells.Item(2, 14).Value = 'Cumple'
$ws.Cells.Item(2, 15).Value = 'No aplica'
$ws.Cells.Item(2, 16).Value = 'Cumple'
$ws.Cells.Item(2, 17).Value = 'VERA BURNEO LUIS RAMIRO'
$ws.Cells.Item(2, 18).Value = 'FJSLKJDASLKJDA'

# Row 3
$ws.Cells.Item(3, 1).Value = '2024-08-29'
$ws.Cells.Item(3, 2).Value = 'II'
$ws.Cells.Item(3, 3).Value = 'dsalkdjlsakd'
$ws.Cells.Item(3, 4).Value = 'BONE ERAZO SANDRA ELIZABETH'
$ws.Cells.Item(3, 5).Value = 'Cumple'
$ws.Cells.Item(3, 6).Value = 'No cumple'
$ws.Cells.Item(3, 7).Value = 'No cumple'
$ws.Cells.Item(3, 8).Value = 'No cumple'
$ws.Cells.Item(3, 9).Value = 'Cumple'
$ws.Cells.Item(3, 10).Value = 'No cumple'
$ws.Cells.Item(3, 11).Value = 'No cumple'
$ws.Cells.Item(3, 12).Value = 'Cumple'
$ws.Cells.Item(3, 13).Value = 'Cumple'
$ws.Cells.Item(3, 14).Value = 'Cumple'
$ws.Cells.Item(3, 15).Value = 'Cumple'
$ws.Cells.Item(3, 16).Value = 'Cumple'
$ws.Cells.Item(3, 17).Value = 'GAVILANEZ QUISPE EDWIN ORLANDO'
$ws.Cells.Item(3, 18).Value = 'dsadddddddddddddddddddddddd'

# Row 4
$ws.Cells.Item(4, 1).Value = '2024-08-29'
$ws.Cells.Item(4, 2).Value = 'II'
$ws.Cells.Item(4, 3).Value = 'dsadsada'
$ws.Cells.Item(4, 4).Value = 'BALSECA ALEGRIA MARCELA DEL PILAR'
$ws.Cells.Item(4, 5).Value = 'No cumple'
$ws.Cells.Item(4, 6).Value = 'No cumple'
$ws.Cells.Item(4, 7).Value = 'Cumple'
$ws.Cells.Item(4, 8).Value = 'No aplica'
$ws.Cells.Item(4, 9).Value = 'Cumple'
$ws.Cells.Item(4, 10).Value = 'No cumple'
$ws.Cells.Item(4, 11).Value = 'Cumple'
$ws.Cells.Item(4, 12).Value = 'Cumple'
$ws.Cells.Item(4, 13).Value = 'Cumple'
$ws.Cells.Item(4, 14).Value = 'Cumple'
$ws.Cells.Item(4, 15).Value = 'Cumple'
$ws.Cells.Item(4, 16).Value = 'Cumple'
$ws.Cells.Item(4, 17).Value = 'VILLA BALCAZAR LENIN ARMANDO'
$ws.Cells.Item(4, 18).Value = 'dasdsad'

# Row 5
$ws.Cells.Item(5, 1).Value = '2024-08-29'
$ws.Cells.Item(5, 2).Value = 'II'
$ws.Cells.Item(5, 3).Value = 'dsadsada'
$ws.Cells.Item(5, 4).Value = 'BOLANOS ORTIZ EDGAR ARTURO'
$ws.Cells.Item(5, 5).Value = 'Cumple'
$ws.Cells.Item(5, 17).Value = 'RIOS ALCIVAR ERICK RICARDO'
$ws.Cells.Item(5, 18).Value = 'sdgsgfd'

# Row 6
$ws.Cells.Item(6, 1).Value = '2024-08-30'
$ws.Cells.Item(6, 2).Value = 'III'
$ws.Cells.Item(6, 3).Value = 'BODEGA 1'
$ws.Cells.Item(6, 4).Value = 'BALSECA ALEGRIA MARCELA DEL PILAR'
$ws.Cells.Item(6, 5).Value = 'Cumple'
$ws.Cells.Item(6, 6).Value = 'Cumple'
$ws.Cells.Item(6, 7).Value = 'Cumple'
$ws.Cells.Item(6, 8).Value = 'No cumple'
$ws.Cells.Item(6, 9).Value = 'No cumple'
$ws.Cells.Item(6, 10).Value = 'Cumple'
$ws.Cells.Item(6, 11).Value = 'Cumple'
$ws.Cells.Item(6, 12).Value = 'Cumple'
$ws.Cells.Item(6, 13).Value = 'No cumple'
$ws.Cells.Item(6, 14).Value = 'No cumple'
$ws.Cells.Item(6, 15).Value = 'No cumple'
$ws.Cells.Item(6, 16).Value = 'No cumple'
$ws.Cells.Item(6, 17).Value = 'QUEZADA ALBAN DARWIN EDUARDO'
$ws.Cells.Item(6, 18).Value = 'sadadadad'

# Row 7
$ws.Cells.Item(7, 1).Value = '2024-08-30'
$ws.Cells.Item(7, 2).Value = 'I'
$ws.Cells.Item(7, 3).Value = 'BODEGA 5 - PASTIFICIO'
$ws.Cells.Item(7, 4).Value = 'CORREA ESPINOSA DIANA JESSELA'
$ws.Cells.Item(7, 5).Value = 'No cumple'
$ws.Cells.Item(7, 6).Value = 'Cumple'
$ws.Cells.Item(7, 7).Value = 'No cumple'
$ws.Cells.Item(7, 8).Value = 'No cumple'
$ws.Cells.Item(7, 9).Value = 'Cumple'
$ws.Cells.Item(7, 10).Value = 'No cumple'
$ws.Cells.Item(7, 11).Value = 'No cumple'
$ws.Cells.Item(7, 12).Value = 'No cumple'
$ws.Cells.Item(7, 13).Value = 'Cumple'
$ws.Cells.Item(7, 14).Value = 'No cumple'
$ws.Cells.Item(7, 15).Value = 'No cumple'
$ws.Cells.Item(7, 16).Value = 'Cumple'
$ws.Cells.Item(7, 17).Value = 'QUISPE TOAPANTA SEGUNDO ARMANDO'
$ws.Cells.Item(7, 18).Value = 'tttttttttttttttttt'

# Row 5 columns F:P have no value in the updated data (cleared).
$ws.Range("F5:P5").ClearContents()

# The former last row (row 8) is removed entirely.
$ws.Rows.Item(8).Delete()
